$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove cells that no longer have values in the target
$ws.Range("E2").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("AO4").ClearContents()

# Update cell values (numbers and company-name text) per the capital structure refresh
$ws.Range("D2").Value = 0.42505
$ws.Range("G2").Value = 0.4893931524909113
$ws.Range("H2").Value = 0.4893931524909113
$ws.Range("I2").Value = 0.7681514921497343
$ws.Range("J2").Value = 0.7069359527926456
$ws.Range("K2").Value = -104.789
$ws.Range("L2").Value = -0.4186368902560825
$ws.Range("M2").Value = 44.29724
$ws.Range("N2").Value = 0.02032999986231648
$ws.Range("O2").Value = -0.4227279580872038
$ws.Range("P2").Value = 44.29724
$ws.Range("Q2").Value = 0.02032999986231648
$ws.Range("R2").Value = -0.4227279580872038
$ws.Range("U2").Value = 420.15
$ws.Range("V2").Value = 0.1928257706835069
$ws.Range("W2").Value = 0.09446718948218198
$ws.Range("X2").Value = 0.03928648001874933
$ws.Range("Y2").Value = 0.05518070946343265
$ws.Range("Z2").Value = 0.1581424969342556
$ws.Range("AA2").Value = 0.1438616301288109
$ws.Range("AB2").Value = 0.03846338425444983
$ws.Range("AC2").Value = 0.1054243970196531
$ws.Range("AD2").Value = 297.384
$ws.Range("AF2").Value = 297.384
$ws.Range("AG2").Value = -122.766
$ws.Range("AH2").Value = 0.1200923638307891
$ws.Range("AI2").Value = 0.1346833186595798
$ws.Range("AJ2").Value = -0.05970690768739931
$ws.Range("AK2").Value = -0.06866591269854586
$ws.Range("AL2").Value = 10.699
$ws.Range("AM2").Value = 9.708
$ws.Range("AN2").Value = 1.879678907780798
$ws.Range("AO2").Value = 17.97139919618656
$ws.Range("AP2").Value = -0.7759686492636366
$ws.Range("AQ2").Value = 19.80593325092707
$ws.Range("B3").Value = "Beluga NV (ENXTBR:BELU)"
$ws.Range("D3").Value = 0.9420000000000001
$ws.Range("G3").Value = 0.4054054054054054
$ws.Range("H3").Value = 0.4054054054054054
$ws.Range("I3").Value = 0.3633633633633633
$ws.Range("J3").Value = 0.3456493993993994
$ws.Range("K3").Value = 2.3
$ws.Range("L3").Value = 0.3453453453453453
$ws.Range("M3").Value = 0.53176
$ws.Range("N3").Value = 0.1100952380952381
$ws.Range("O3").Value = 0.2312
$ws.Range("P3").Value = 0.53176
$ws.Range("Q3").Value = 0.1100952380952381
$ws.Range("R3").Value = 0.2312
$ws.Range("U3").Value = 4.68
$ws.Range("V3").Value = 0.9689440993788819
$ws.Range("W3").Value = 0.6534090909090908
$ws.Range("X3").Value = 0.04010957578304883
$ws.Range("Y3").Value = 0.613299515126042
$ws.Range("Z3").Value = 2.917214191852825
$ws.Range("AA3").Value = 1.008333333333333
$ws.Range("AB3").Value = 0.0390025566093618
$ws.Range("AC3").Value = 0.9693307767239713
$ws.Range("AD3").Value = 0.384
$ws.Range("AF3").Value = 0.384
$ws.Range("AG3").Value = -4.295999999999999
$ws.Range("AH3").Value = 0.07364787111622555
$ws.Range("AI3").Value = 0.06260189109879361
$ws.Range("AJ3").Value = -8.044943820224708
$ws.Range("AK3").Value = -2.954607977991745
$ws.Range("AL3").Value = 0.065
$ws.Range("AM3").Value = 0.017
$ws.Range("AN3").Value = 0.1422222222222222
$ws.Range("AO3").Value = 37.23076923076923
$ws.Range("AP3").Value = -1.591111111111111
$ws.Range("AQ3").Value = 142.3529411764706
$ws.Range("G4").Value = 0.1104972375690608
$ws.Range("H4").Value = 0.1104972375690608
$ws.Range("I4").Value = 0.08353591160220994
$ws.Range("J4").Value = 0.06414927551339518
$ws.Range("K4").Value = 0.8110000000000001
$ws.Range("L4").Value = 0.08961325966850829
$ws.Range("M4").Value = 0.7354800000000001
$ws.Range("N4").Value = 0.1476867469879518
$ws.Range("O4").Value = 0.9068803945745993
$ws.Range("P4").Value = 0.7354800000000001
$ws.Range("Q4").Value = 0.1476867469879518
$ws.Range("R4").Value = 0.9068803945745993
$ws.Range("U4").Value = 3.57
$ws.Range("V4").Value = 0.716867469879518
$ws.Range("W4").Value = 0.13561872909699
$ws.Range("X4").Value = 0.03846338425444983
$ws.Range("Y4").Value = 0.09715534484254014
$ws.Range("Z4").Value = 2.563739376770537
$ws.Range("AA4").Value = 0.1644620236249933
$ws.Range("AB4").Value = 0.03846338425444983
$ws.Range("AC4").Value = 0.1259986393705435
$ws.Range("AD4").Value = 0
$ws.Range("AF4").Value = 0
$ws.Range("AG4").Value = -3.57
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = -2.53191489361702
$ws.Range("AK4").Value = -1.475206611570248
$ws.Range("AL4").Value = 0
$ws.Range("AM4").Value = -0.478
$ws.Range("AN4").Value = 0
$ws.Range("AP4").Value = -3.216216216216216
$ws.Range("AQ4").Value = -1.581589958158996
$ws.Range("B5").Value = "GIMV NV (ENXTBR:GIMB)"
$ws.Range("D5").Value = -0.0919
$ws.Range("G5").Value = 0.6107969151670951
$ws.Range("H5").Value = 0.6107969151670951
$ws.Range("I5").Value = 0.7881748071979435
$ws.Range("J5").Value = 0.7881748071979435
$ws.Range("K5").Value = -128
$ws.Range("L5").Value = -0.6580976863753213
$ws.Range("M5").Value = 41.5
$ws.Range("N5").Value = 0.02604983993471847
$ws.Range("O5").Value = -0.32421875
$ws.Range("P5").Value = 41.5
$ws.Range("Q5").Value = 0.02604983993471847
$ws.Range("R5").Value = -0.32421875
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 355.7
$ws.Range("V5").Value = 0.2232753750549244
$ws.Range("W5").Value = -0.08976786590925029
$ws.Range("X5").Value = 0.04232358313571015
$ws.Range("Y5").Value = -0.1320914490449605
$ws.Range("Z5").Value = 0.1563881965104125
$ws.Range("AA5").Value = 0.1232612366326284
$ws.Range("AB5").Value = 0.03841108196386579
$ws.Range("AC5").Value = 0.08485015466876265
$ws.Range("AD5").Value = 297
$ws.Range("AF5").Value = 297
$ws.Range("AG5").Value = -58.69999999999999
$ws.Range("AH5").Value = 0.1571345431458653
$ws.Range("AI5").Value = 0.1752936315882666
$ws.Range("AJ5").Value = -0.03825599582898852
$ws.Range("AK5").Value = -0.04385178544748244
$ws.Range("AL5").Value = 10.5
$ws.Range("AM5").Value = 10.095
$ws.Range("AN5").Value = 1.923575129533679
$ws.Range("AO5").Value = 14.6
$ws.Range("AP5").Value = -0.380181347150259
$ws.Range("AQ5").Value = 15.18573551263001
$ws.Range("I6").Value = 0.8927680798004987
$ws.Range("J6").Value = 0.8588941521793602
$ws.Range("K6").Value = 20.1
$ws.Range("L6").Value = 0.5012468827930174
$ws.Range("M6").Value = 1.53
$ws.Range("N6").Value = 0.00265625
$ws.Range("O6").Value = 0.07611940298507462
$ws.Range("P6").Value = 1.53
$ws.Range("Q6").Value = 0.00265625
$ws.Range("R6").Value = 0.07611940298507462
$ws.Range("U6").Value = 56.2
$ws.Range("V6").Value = 0.09756944444444444
$ws.Range("W6").Value = 0.05331564986737401
$ws.Range("X6").Value = 0.03846338425444983
$ws.Range("Y6").Value = 0.01485226561292418
$ws.Range("Z6").Value = 0.1203120312031203
$ws.Range("AA6").Value = 0.1033353000371807
$ws.Range("AB6").Value = 0.03846338425444983
$ws.Range("AC6").Value = 0.06487191578273092
$ws.Range("AG6").Value = -56.2
$ws.Range("AJ6").Value = -0.1081185071181224
$ws.Range("AK6").Value = -0.1261787157611136
$ws.Range("AL6").Value = 0.134
$ws.Range("AM6").Value = 0.07400000000000001
$ws.Range("AO6").Value = 267.1641791044776
$ws.Range("AQ6").Value = 483.7837837837837
